$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 329
$ws.Range("F3").Value = 93
$ws.Range("F4").Value = 15
$ws.Range("F5").Value = 324
$ws.Range("F6").Value = 82
$ws.Range("F7").Value = 50
$ws.Range("F8").Value = 133
$ws.Range("F11").Value = 473
$ws.Range("F12").Value = 271
$ws.Range("F13").Value = 105
$ws.Range("F15").Value = 86
